# Template standardization / file load for hipot changes
# - Enable the "Enable" flag (column E) for the CHECK_INTERLOCK (row 6)
#   and HIPOT_TEST (row 7) steps in the Test Flow sheet.
# - Update the active selection to the single cell E8 (matching the
#   active cell) instead of the previous multi-cell range A6:E8.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E6").Value = $true
$ws.Range("E7").Value = $true

$ws.Range("E8").Select() | Out-Null
